# Insert a new data row above row 16 (pushes existing rows 16-86 down to 17-87)
# so the sheet grows from A1:R86 to A1:R87, matching the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(16).Insert()

# Populate the newly inserted row 16 with its full record (mirrors the
# surrounding rows' static columns, with the new unique measurement values).
$ws.Range("A16").Value = 11
$ws.Range("B16").Value = "Vega Monumental Concepción"
$ws.Range("C16").Value = "Bíobío"
$ws.Range("D16").Value = 44687
$ws.Range("E16").Value = 8
$ws.Range("F16").Value = 100112001
$ws.Range("G16").Value = "Berenjena"
$ws.Range("H16").Value = "Sin especificar"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 150
$ws.Range("K16").Value = 5500
$ws.Range("L16").Value = 6000
$ws.Range("M16").Value = 5767
$ws.Range("N16").Value = "$/caja 60 unidades"
$ws.Range("O16").Value = "Región de Arica y Parinacota"
$ws.Range("P16").Value = 96
$ws.Range("Q16").Value = 60
$ws.Range("R16").Value = "Hortaliza"
